$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Cost Of Investment *" header in E1 with the new
# "Face Value For Redemption *" header (this also removes the now-unused
# shared string and appends the new one, shifting other shared string
# indices down automatically).
$ws.Range("E1").Value = "Face Value For Redemption *"

# Widen column E to fit the new, longer header text.
$ws.Columns("E:E").ColumnWidth = 28.285714285714285

# Update the active selection/cell to E2.
[void]$ws.Range("E2").Select()

Write-Host "done"
